$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 and Row 3 swap company/ticker labels (shared-string reorder effect) ---
$ws.Range("B2").Value = "HANWHA AEROSPACE"
$ws.Range("C2").Value = "012450.KS"
$ws.Range("B3").Value = "HYUNDAI ROTEM"
$ws.Range("C3").Value = "064350.KS"

# --- Row 2 numeric updates ---
$ws.Range("D2").Value = 875000
$ws.Range("E2").Value = 32.9
$ws.Range("F2").Value = 2.7
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 54.7
$ws.Range("N2").Value = 54.84087454262382

# --- Row 3 numeric updates ---
$ws.Range("D3").Value = 179700
$ws.Range("E3").Value = 33.8
$ws.Range("F3").Value = 2.28
$ws.Range("G3").Value = 30
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 51.9
$ws.Range("N3").Value = 54.84087454262382

# --- Row 4 numeric updates ---
$ws.Range("D4").Value = 46800
$ws.Range("E4").Value = 23.2
$ws.Range("F4").Value = 1.3
$ws.Range("N4").Value = 54.84087454262382

# --- Row 5 numeric updates ---
$ws.Range("D5").Value = 104900
$ws.Range("E5").Value = 38.1
$ws.Range("F5").Value = -3.67
$ws.Range("N5").Value = 54.84087454262382

# --- Row 6 numeric updates ---
$ws.Range("D6").Value = 366000
$ws.Range("E6").Value = 26.6
$ws.Range("F6").Value = -4.44
$ws.Range("N6").Value = 54.84087454262382
